$d = $word.ActiveDocument

function Get-Text {
    return $d.Content.Text
}

function Replace-Range([int]$start, [int]$len, [string]$newText) {
    $r = $d.Range($start, $start + $len)
    $r.Text = $newText
}

# Replace exactly the (unique) substring $old with $new - used only when
# $old lies completely inside a single run (or is made up of whole runs
# we *want* merged) so we get the formatting behaviour we expect.
function Replace-InRun([string]$old, [string]$new) {
    $t = Get-Text
    $start = $t.IndexOf($old)
    if ($start -lt 0) {
        throw "anchor not found: $old"
    }
    Replace-Range $start $old.Length $new
}

# Delete the (unique) substring $old - used when $old is made up of whole
# runs so deleting it just drops those runs and leaves the neighbours
# (and their formatting) untouched.
function Delete-Run([string]$old) {
    $t = Get-Text
    $start = $t.IndexOf($old)
    if ($start -lt 0) {
        throw "anchor not found: $old"
    }
    $r = $d.Range($start, $start + $old.Length)
    $r.Delete()
}

# Locate the unique $prefix (used only for positioning, never touched,
# may itself span several runs) and replace the $len characters that
# immediately follow it. Used to edit a run in place without disturbing
# the run(s) that precede it, even when the run's own text (e.g. a lone
# space) is too short/common to be a unique anchor by itself.
function Replace-After([string]$prefix, [int]$len, [string]$newText) {
    $t = Get-Text
    $pstart = $t.IndexOf($prefix)
    if ($pstart -lt 0) {
        throw "anchor not found: $prefix"
    }
    $start = $pstart + $prefix.Length
    Replace-Range $start $len $newText
}

# -----------------------------------------------------------------
# 1) "<id>p050v_1</id>" and "<id>p050v_2</id>" : collapse the three
#    runs (<id>, text, </id>) into a single run. Setting identical
#    text is a no-op in this engine, so swap through a placeholder to
#    force the merge; the merged run keeps the first run's formatting,
#    which is what the target wants here.
# -----------------------------------------------------------------
$ph = [char]0xE000
$placeholder1 = "$ph" + "PLACEHOLDER1" + "$ph"
$placeholder2 = "$ph" + "PLACEHOLDER2" + "$ph"

Replace-InRun "<id>p050v_1</id>" $placeholder1
Replace-InRun $placeholder1 "<id>p050v_1</id>"

Replace-InRun "<id>p050v_2</id>" $placeholder2
Replace-InRun $placeholder2 "<id>p050v_2</id>"

# -----------------------------------------------------------------
# 2) "... with <del><fr>de la <m>mie</m></fr></del> <m><add>very small
#    crumbs of the inside of a loaf</add> of bread</m> and ..."
#    becomes
#    "... with <del><m>pith</m></del> <m><add>bread pith</add></m> and ..."
#
#    Run layout (before):
#      <del> | <fr> | "de " | "la " | <m> | "mie" | </m></fr> | </del> | " "
#      | <m> | <add> | "very " | "small" | " " | "crumbs of the inside of a loaf"
#      | </add> | " of bread" | </m> | " and "
# -----------------------------------------------------------------

# 2a) "mie" -> "pith", touching only that run (anchor on the preceding,
#     still-intact "<fr>de la <m>" text, none of which gets replaced)
Replace-After "<fr>de la <m>" 3 "pith"

# 2b) "</m></fr>" -> "</m>", touching only that run
Replace-After "<fr>de la <m>pith" 9 "</m>"

# 2c) now drop the "<fr>" + "de " + "la " runs entirely (whole runs only)
Delete-Run "<fr>de la "

# 2d) the lone space run between </del> and the next <m> gains text;
#     anchor on the preceding (unique, untouched) "</del>" run and only
#     replace the 1 character that belongs to the space run itself
Replace-After "</del>" 1 " a very little amount of "

# 2e) drop the "very " + "small" + " " runs, leaving the
#     "crumbs of the inside of a loaf" run (and its own formatting) alone
Delete-Run "very small "

# 2f) rename that surviving run's text - stays fully inside its own run
Replace-InRun "crumbs of the inside of a loaf" "bread pith"

# 2g) drop the trailing " of bread" run entirely
Delete-Run " of bread"
